$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data table (columns A:E, rows 4:11) — the "Step3Fields" entry is removed
# and two new entries ("FieldMappingsCardUI" / "ImportSettingsCardUI") are
# added, so every row from the old "ParentAttribute" row downward shifts.
# The summary box living in columns F:I (rows 1,2,4,5,6) is NOT part of this
# table and must stay untouched, so we write cell-by-cell instead of using a
# structural row insert/delete (which would also drag the F:I box and the
# far-below formatting-only rows along with it).
# ---------------------------------------------------------------------------

$rows = @(
  @{ Row = 4;  A = "FieldMappingsCardUI";   B = 7; C = 7;    D = "Automated";         E = $null },
  @{ Row = 5;  A = "ImportSettingsCardUI";  B = 4; C = 4;    D = "Automated";         E = $null },
  @{ Row = 6;  A = "ParentAttribute";       B = 0; C = $null; D = "Unwritten";        E = $null },
  @{ Row = 7;  A = "RDOSpecificFields";     B = 1; C = 1;    D = "Automated";         E = "Long run time. (~1h)" },
  @{ Row = 8;  A = "SavingErrors";          B = 0; C = 9;    D = "Ready to Write";    E = $null },
  @{ Row = 9;  A = "SourceAttributeFields"; B = 0; C = 3;    D = "Suited to Manual";  E = $null },
  @{ Row = 10; A = "SourceAttributeLists";  B = 0; C = 11;   D = "Suited to Manual";  E = $null },
  @{ Row = 11; A = "WorkspaceFields";       B = 5; C = 13;   D = "Ready to Write";    E = $null }
)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("A$n").Value2 = $r.A
  $ws.Range("B$n").Value2 = $r.B

  if ($null -eq $r.C) {
    $ws.Range("C$n").Clear()
  } else {
    $ws.Range("C$n").Value2 = $r.C
  }

  $ws.Range("D$n").Value2 = $r.D

  if ($null -eq $r.E) {
    $ws.Range("E$n").Clear()
  } else {
    $ws.Range("E$n").Value2 = $r.E
  }
}

# ---------------------------------------------------------------------------
# Summary formulas in the fixed F:I box — their ranges grow by one row
# (one net extra data row: +2 new, -1 removed) so update the formula text.
# ---------------------------------------------------------------------------
$ws.Range("E1").Formula = "=COUNTA(`$A`$3:A40)"
$ws.Range("G1").Formula = '=COUNTIF($D$3:D39,"Ready to Write")+COUNTIF($D$3:D39,"Outdated")+COUNTIF($D$3:D39,"Writing")+COUNTIF($D$3:D39,"Testing")'
$ws.Range("G2").Formula = '=COUNTIF($D$3:D39,"Automated")+COUNTIF($D$3:D39,"Finished")'
$ws.Range("G4").Formula = "=SUM(`$C`$3:C39)"
$ws.Range("G5").Formula = "=SUM(`$B`$3:B39)"
$ws.Range("G6").Formula = "=G5/G4"

# ---------------------------------------------------------------------------
# Conditional formatting — the "no blanks" rule now spans one row further.
# ---------------------------------------------------------------------------
$ws.Range("D2:D52").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D53"))

# ---------------------------------------------------------------------------
# Sheet view — scrolled down with a new active selection.
# ---------------------------------------------------------------------------
$ws.Range("D5").Select()
$excel.ActiveWindow.ScrollRow = 22
